# Fruta / hortaliza, semanal
# Insert a new week of price data (2021-09-29, serial 44468) for
# Comercializadora del Agro de Limarí - Chirimoya, ahead of the existing
# rows, shifting the existing data block down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows starting at row 24 (pushes old rows 24-32 down to 27-35)
$ws.Rows("24:26").Insert()

# Common / repeated column values for this product block
$colA = 2
$colB = "Comercializadora del Agro de Limarí"
$colC = "Coquimbo"
$colE = 4
$colF = "Fruta"
$colG = 100107
$colH = "Otros"
$colI = 100107002
$colJ = "Chirimoya"
$colK = "Cultivar IV Región"
$colQ = "$/kilo (en caja de 15 kilos)"
$colR = "Provincia de Limarí"

# New week date (serial 44468 = 2021-09-29)
$newDate = 44468

# Row 24 - Calidad "Especial"
$r = 24
$ws.Cells.Item($r, 1).Value = $colA
$ws.Cells.Item($r, 2).Value = $colB
$ws.Cells.Item($r, 3).Value = $colC
$ws.Cells.Item($r, 4).Value = $newDate
$ws.Cells.Item($r, 5).Value = $colE
$ws.Cells.Item($r, 6).Value = $colF
$ws.Cells.Item($r, 7).Value = $colG
$ws.Cells.Item($r, 8).Value = $colH
$ws.Cells.Item($r, 9).Value = $colI
$ws.Cells.Item($r, 10).Value = $colJ
$ws.Cells.Item($r, 11).Value = $colK
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 340
$ws.Cells.Item($r, 14).Value = 2500
$ws.Cells.Item($r, 15).Value = 2600
$ws.Cells.Item($r, 16).Value = 2550
$ws.Cells.Item($r, 17).Value = $colQ
$ws.Cells.Item($r, 18).Value = $colR
$ws.Cells.Item($r, 19).Value = 2550
$ws.Cells.Item($r, 20).Value = 1

# Row 25 - Calidad "Primera"
$r = 25
$ws.Cells.Item($r, 1).Value = $colA
$ws.Cells.Item($r, 2).Value = $colB
$ws.Cells.Item($r, 3).Value = $colC
$ws.Cells.Item($r, 4).Value = $newDate
$ws.Cells.Item($r, 5).Value = $colE
$ws.Cells.Item($r, 6).Value = $colF
$ws.Cells.Item($r, 7).Value = $colG
$ws.Cells.Item($r, 8).Value = $colH
$ws.Cells.Item($r, 9).Value = $colI
$ws.Cells.Item($r, 10).Value = $colJ
$ws.Cells.Item($r, 11).Value = $colK
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 400
$ws.Cells.Item($r, 14).Value = 2100
$ws.Cells.Item($r, 15).Value = 2200
$ws.Cells.Item($r, 16).Value = 2150
$ws.Cells.Item($r, 17).Value = $colQ
$ws.Cells.Item($r, 18).Value = $colR
$ws.Cells.Item($r, 19).Value = 2150
$ws.Cells.Item($r, 20).Value = 1

# Row 26 - Calidad "Segunda"
$r = 26
$ws.Cells.Item($r, 1).Value = $colA
$ws.Cells.Item($r, 2).Value = $colB
$ws.Cells.Item($r, 3).Value = $colC
$ws.Cells.Item($r, 4).Value = $newDate
$ws.Cells.Item($r, 5).Value = $colE
$ws.Cells.Item($r, 6).Value = $colF
$ws.Cells.Item($r, 7).Value = $colG
$ws.Cells.Item($r, 8).Value = $colH
$ws.Cells.Item($r, 9).Value = $colI
$ws.Cells.Item($r, 10).Value = $colJ
$ws.Cells.Item($r, 11).Value = $colK
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 240
$ws.Cells.Item($r, 14).Value = 1700
$ws.Cells.Item($r, 15).Value = 1800
$ws.Cells.Item($r, 16).Value = 1750
$ws.Cells.Item($r, 17).Value = $colQ
$ws.Cells.Item($r, 18).Value = $colR
$ws.Cells.Item($r, 19).Value = 1750
$ws.Cells.Item($r, 20).Value = 1

# Keep the sheet's used-range/dimension honest
$ws.Cells.Item(1,1).Select()
